$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new day of data (2023-12-09) below the existing weather/rain log.
$ws.Range("A70").Value = 20231209
$ws.Range("B70").Value = 0

# Match the vertical-center alignment style used by the rest of column B
# (same cell style as B69, i.e. style index "1" / vertical center alignment).
$ws.Range("B70").VerticalAlignment = -4108

# Reflect the updated scroll/selection position recorded for the sheet view.
$ws.Range("I67").Select()

# The workbook's default font was re-baked from "新細明體" to "Calibri"
# (e.g. after being opened/re-saved on a machine without that font
# installed). Update the workbook's Normal ("一般") style font to match.
$normalStyle = $wb.Styles.Item("一般")
$normalStyle.Font.Name = "Calibri"
